$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tbl = $sh.Table

# Header row (row 1): Col0..Col4
for ($c = 1; $c -le 5; $c++) {
    $tbl.Cell(1, $c).Shape.TextFrame.TextRange.Text = "Col" + ($c - 1)
}

# Data rows 2..9 (0-based index 0..7): column 1 -> {{cell0.N}}, columns 3..5 -> {{cellC.N}}
for ($r = 2; $r -le 9; $r++) {
    $n = $r - 2
    $tbl.Cell($r, 1).Shape.TextFrame.TextRange.Text = "{{cell0." + $n + "}}"
    $tbl.Cell($r, 3).Shape.TextFrame.TextRange.Text = "{{cell2." + $n + "}}"
    $tbl.Cell($r, 4).Shape.TextFrame.TextRange.Text = "{{cell3." + $n + "}}"
    $tbl.Cell($r, 5).Shape.TextFrame.TextRange.Text = "{{cell4." + $n + "}}"
}

# Last data row (row 10): clear column 1 text (leave empty)
$tbl.Cell(10, 1).Shape.TextFrame.TextRange.Text = ""
